$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-patch) table contents, already sorted by Mod then Tier
# (mirrors the Mod Part Patching for 0.8.0 commit: fills in missing Tier
# values, renames "spaceplanes" -> "aerodynamics" and "ntr" -> "reactors",
# and re-sorts the data rows by Mod/Tier).
$rows = @(
    @{ Part = "RREngineNTJ1";        Mod = "aerodynamics"; Tier = 11 },
    @{ Part = "RRRCSLinearCold";     Mod = "control";       Tier = 9 },
    @{ Part = "RRRCSQuadCold";       Mod = "control";       Tier = 9 },
    @{ Part = "RREngineMAPT0";       Mod = "ion";           Tier = 10 },
    @{ Part = "RREngineMET";         Mod = "ion";           Tier = 10 },
    @{ Part = "rr.boxblue";          Mod = "isru";          Tier = 9 },
    @{ Part = "rr.boxblue.flat";     Mod = "isru";          Tier = 9 },
    @{ Part = "rr.boxred.flat";      Mod = "isru";          Tier = 9 },
    @{ Part = "rr.boxblue.long";     Mod = "isru";          Tier = 10 },
    @{ Part = "rr.exoscoop";         Mod = "isru";          Tier = 10 },
    @{ Part = "rr.hydroscoop.temp";  Mod = "isru";          Tier = 10 },
    @{ Part = "RREngineNERVA-R";     Mod = "reactors";      Tier = 8 },
    @{ Part = "RREngineNERVA-O";     Mod = "reactors";      Tier = 8 },
    @{ Part = "RREngineSFRJ";        Mod = "solids";        Tier = 9 },
    @{ Part = "RREngineSRB.RT05";    Mod = "solids";        Tier = 10 },
    @{ Part = "RREngineSRB.RT10";    Mod = "solids";        Tier = 10 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Part
    $ws.Cells.Item($r, 2).Value = $row.Mod
    $ws.Cells.Item($r, 3).Value = $row.Tier
    # Column D (Mod dependency) is unchanged for every row.
    $ws.Cells.Item($r, 4).Value = "RationalResourcesParts"

    $formula = '="@PART["&A' + $r + '&"]:AFTER["&D' + $r + '&"] //' + "`n{`n`t@TechRequired = " + '"&B' + $r + '&C' + $r + '&"' + "`n}`""
    $ws.Cells.Item($r, 5).Formula = $formula

    $r++
}

# Setting a multi-line formula triggers an implicit row auto-height; put
# the rows back to the sheet's default (unsized) state like the source file.
$ws.Rows("2:17").EntireRow.AutoFit() | Out-Null

# Matches the saved selection recorded in the commit.
$ws.Range("E16:E17").Select()
